# Quarterly income-statement refresh: drop the oldest quarter (column D),
# shift everything one column to the left, and append the newest quarter
# of data in the now-empty last column (M). Also correct the publish-date
# label for the column that was re-published with a later revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop the oldest quarter column (D). This shifts D:M -> drops D,
#    and slides old E:M left into D:L (carrying values + column widths).
$ws.Columns("D:D").Delete()

# 2. The freshly-vacated last column (M) needs its width restored to match
#    the "wide" column pattern used for this quarter's position.
$ws.Columns("M:M").ColumnWidth = 30.166666666666668

# 3. Populate the new quarter column (M) with the latest data.
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-30"

$ws.Range("M11").Value = 12379
$ws.Range("M12").Value = -12144
$ws.Range("M13").Value = 234
$ws.Range("M14").Value = -202
$ws.Range("M15").Value = "-"
$ws.Range("M16").Value = 7
$ws.Range("M17").Value = 39
$ws.Range("M18").Value = -23
$ws.Range("M19").Value = 302
$ws.Range("M20").Value = 318
$ws.Range("M21").Value = 135
$ws.Range("M22").Value = 453
$ws.Range("M23").Value = "-"
$ws.Range("M24").Value = 453
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 428
$ws.Range("M27").Value = 0

# 4. Column I's quarter (Q4 1400/12) was re-published with a later revision
#    by the time of this refresh; correct its publish-date label.
$ws.Range("I9").Value = "1402-02-30 (8)"
